$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Sending cluster / Target cluster reassigned to ECs; recomputed metrics) ---
# --- Add new rows 3-7 covering ECs/FAPs/sCs sending clusters x ECs/sCs target clusters for Gm13306-Ccr10 ---

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gm13306"
$ws.Cells.Item(2,3).Value = "Ccr10"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.041225
$ws.Cells.Item(2,8).Value = 0.123675
$ws.Cells.Item(2,9).Value = 0.04646624533012825
$ws.Cells.Item(2,10).Value = 0.05129229203448434
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.5706193333333334
$ws.Cells.Item(2,14).Value = 1.711858
$ws.Cells.Item(2,15).Value = 0.2443278446591134
$ws.Cells.Item(2,16).Value = 0.3265937887468804
$ws.Cells.Item(2,17).Value = 0.02352378201666667
$ws.Cells.Item(2,18).Value = 0.21171403815
$ws.Cells.Item(2,19).Value = 0.01135299757091183
$ws.Cells.Item(2,20).Value = 0.01675174398905368

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gm13306"
$ws.Cells.Item(3,3).Value = "Ccr10"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.041225
$ws.Cells.Item(3,8).Value = 0.123675
$ws.Cells.Item(3,9).Value = 0.04646624533012825
$ws.Cells.Item(3,10).Value = 0.05129229203448434
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.7648465
$ws.Cells.Item(3,14).Value = 3.529693
$ws.Cells.Item(3,15).Value = 0.7556721553408867
$ws.Cells.Item(3,16).Value = 0.6734062112531195
$ws.Cells.Item(3,17).Value = 0.0727557969625
$ws.Cells.Item(3,18).Value = 0.436534781775
$ws.Cells.Item(3,19).Value = 0.03511324775921643
$ws.Cells.Item(3,20).Value = 0.03454054804543066

$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Gm13306"
$ws.Cells.Item(4,3).Value = "Ccr10"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.5955496666666666
$ws.Cells.Item(4,8).Value = 1.786649
$ws.Cells.Item(4,9).Value = 0.6712663897540191
$ws.Cells.Item(4,10).Value = 0.7409850193743232
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.5706193333333334
$ws.Cells.Item(4,14).Value = 1.711858
$ws.Cells.Item(4,15).Value = 0.2443278446591134
$ws.Cells.Item(4,16).Value = 0.3265937887468804
$ws.Cells.Item(4,17).Value = 0.3398321537602222
$ws.Cells.Item(4,18).Value = 3.058489383842
$ws.Cells.Item(4,19).Value = 0.1640090702007038
$ws.Cells.Item(4,20).Value = 0.2420011048821408

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gm13306"
$ws.Cells.Item(5,3).Value = "Ccr10"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.5955496666666666
$ws.Cells.Item(5,8).Value = 1.786649
$ws.Cells.Item(5,9).Value = 0.6712663897540191
$ws.Cells.Item(5,10).Value = 0.7409850193743232
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.7648465
$ws.Cells.Item(5,14).Value = 3.529693
$ws.Cells.Item(5,15).Value = 0.7556721553408867
$ws.Cells.Item(5,16).Value = 0.6734062112531195
$ws.Cells.Item(5,17).Value = 1.051053744792833
$ws.Cells.Item(5,18).Value = 6.306322468757
$ws.Cells.Item(5,19).Value = 0.5072573195533153
$ws.Cells.Item(5,20).Value = 0.4989839144921823

$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Gm13306"
$ws.Cells.Item(6,3).Value = "Ccr10"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.2504285
$ws.Cells.Item(6,8).Value = 0.500857
$ws.Cells.Item(6,9).Value = 0.2822673649158526
$ws.Cells.Item(6,10).Value = 0.2077226885911924
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.5706193333333334
$ws.Cells.Item(6,14).Value = 1.711858
$ws.Cells.Item(6,15).Value = 0.2443278446591134
$ws.Cells.Item(6,16).Value = 0.3265937887468804
$ws.Cells.Item(6,17).Value = 0.1428993437176667
$ws.Cells.Item(6,18).Value = 0.8573960623060001
$ws.Cells.Item(6,19).Value = 0.0689657768874977
$ws.Cells.Item(6,20).Value = 0.06784093987568594

$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Gm13306"
$ws.Cells.Item(7,3).Value = "Ccr10"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.2504285
$ws.Cells.Item(7,8).Value = 0.500857
$ws.Cells.Item(7,9).Value = 0.2822673649158526
$ws.Cells.Item(7,10).Value = 0.2077226885911924
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.7648465
$ws.Cells.Item(7,14).Value = 3.529693
$ws.Cells.Item(7,15).Value = 0.7556721553408867
$ws.Cells.Item(7,16).Value = 0.6734062112531195
$ws.Cells.Item(7,17).Value = 0.44196786172525
$ws.Cells.Item(7,18).Value = 1.767871446901
$ws.Cells.Item(7,19).Value = 0.2133015880283549
$ws.Cells.Item(7,20).Value = 0.1398817487155065

